$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("results")

# --- Header row (row 2): copy alternating style from H2 (s=16) / I2 (s=14), then set the new label ---
$ws.Range("H2").Copy($ws.Range("L2"))
$ws.Range("L2").Value = "5 occ"
$ws.Range("I2").Copy($ws.Range("M2"))
$ws.Range("M2").Value = "6 occ"
$ws.Range("H2").Copy($ws.Range("N2"))
$ws.Range("N2").Value = "7 occ"
$ws.Range("I2").Copy($ws.Range("O2"))
$ws.Range("O2").Value = "8 occ"
$ws.Range("H2").Copy($ws.Range("P2"))
$ws.Range("P2").Value = "9 occ"
$ws.Range("I2").Copy($ws.Range("Q2"))
$ws.Range("Q2").Value = "10 occ"
$ws.Range("H2").Copy($ws.Range("R2"))
$ws.Range("R2").Value = "11 occ"
$ws.Range("I2").Copy($ws.Range("S2"))
$ws.Range("S2").Value = "12 occ"
$ws.Range("H2").Copy($ws.Range("T2"))
$ws.Range("T2").Value = "13 occ"
$ws.Range("I2").Copy($ws.Range("U2"))
$ws.Range("U2").Value = "14 occ"
$ws.Range("H2").Copy($ws.Range("V2"))
$ws.Range("V2").Value = "15 occ"
$ws.Range("I2").Copy($ws.Range("W2"))
$ws.Range("W2").Value = "16 occ"
$ws.Range("H2").Copy($ws.Range("X2"))
$ws.Range("X2").Value = "17 occ"
$ws.Range("I2").Copy($ws.Range("Y2"))
$ws.Range("Y2").Value = "18 occ"
$ws.Range("H2").Copy($ws.Range("Z2"))
$ws.Range("Z2").Value = "19 occ"
$ws.Range("I2").Copy($ws.Range("AA2"))
$ws.Range("AA2").Value = "20 occ"

# --- Row 3 (new unstyled numeric cells) ---
$ws.Range("L3").Value = 8378
$ws.Range("M3").Value = 6545
$ws.Range("N3").Value = 5387
$ws.Range("O3").Value = 4606
$ws.Range("P3").Value = 4079
$ws.Range("Q3").Value = 3827
$ws.Range("R3").Value = 3674
$ws.Range("S3").Value = 3575
$ws.Range("T3").Value = 3515
$ws.Range("U3").Value = 3490
$ws.Range("V3").Value = 3465
$ws.Range("W3").Value = 3414
$ws.Range("X3").Value = 3405
$ws.Range("Y3").Value = 3389
$ws.Range("Z3").Value = 3380
$ws.Range("AA3").Value = 3376

# --- Row 4 (new unstyled numeric cells) ---
$ws.Range("L4").Value = 29211
$ws.Range("M4").Value = 22678
$ws.Range("N4").Value = 10607
$ws.Range("O4").Value = 8504
$ws.Range("P4").Value = 3611
$ws.Range("Q4").Value = 2909
$ws.Range("R4").Value = 778
$ws.Range("S4").Value = 362
$ws.Range("T4").Value = 245
$ws.Range("U4").Value = 239
$ws.Range("V4").Value = 239
$ws.Range("W4").Value = 217
$ws.Range("X4").Value = 210
$ws.Range("Y4").Value = 208
$ws.Range("Z4").Value = 208
$ws.Range("AA4").Value = 96

# --- Row 7: clear styles on F7:K7 (existing cells) and add new unstyled cells L7:AA7 ---
$ws.Range("F7:K7").ClearFormats()
$ws.Range("L7").Value = 1305
$ws.Range("M7").Value = 999
$ws.Range("N7").Value = 892
$ws.Range("O7").Value = 775
$ws.Range("P7").Value = 673
$ws.Range("Q7").Value = 586
$ws.Range("R7").Value = 553
$ws.Range("S7").Value = 476
$ws.Range("T7").Value = 409
$ws.Range("U7").Value = 382
$ws.Range("V7").Value = 376
$ws.Range("W7").Value = 282
$ws.Range("X7").Value = 281
$ws.Range("Y7").Value = 259
$ws.Range("Z7").Value = 149
$ws.Range("AA7").Value = 149

# --- Selection / active sheet ---
$ws.Select()
$ws.Range("I16").Select()
